# Gergely_zeitliste.xlsx: add a new time-tracking entry
#
# Two new rows are inserted above the former row 39 for a CRC research
# session on 2018-06-10 (serial date 43261), pushing every later row
# (and the trailing spacer row) down by two. The SUM() formula in the
# totals row auto-extends to cover the new rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows right before the current row 39 (they inherit
# formatting from the row above, matching style "1" on column A).
[void]$ws.Rows("39:40").Insert()

$activity = "Recherche und Implementierungsversuche CRC"

$ws.Cells.Item(39, 1).Value = 43261
$ws.Cells.Item(39, 2).Value = $activity
$ws.Cells.Item(39, 3).Value = 2

$ws.Cells.Item(40, 1).Value = 43261
$ws.Cells.Item(40, 2).Value = $activity
$ws.Cells.Item(40, 3).Value = 2

# Match the author's final cursor position.
[void]$ws.Range("C39").Select()
